$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (row 1) to make room for headers.
# This shifts the existing data (Monday..Sunday / amounts) down by one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Day"
$ws.Range("B1").Value = "Amount in Liters"

# Move active selection to F4 to mirror the saved view state.
$ws.Range("F4").Select()
